# Update 2018-05-22#1 - Clean up projects: add "Clients" error-code section
# (RegisterClient: 2301/2302/2303) to the db.error.code workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New section header row (merged, bold, left-aligned - same look as the
#     other "table name" header rows already on the sheet) ---------------
$ws.Range("A240").Value = "Clients"
$ws.Range("A234:C234").Copy()
$ws.Range("A240:C240").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A240:C240").Merge()

# --- New error-code rows for the Clients table ---------------------------
$ws.Range("A241").Value = 2301
$ws.Range("B241").Value = "Client Id cannot be null or empty string."
$ws.Range("C241").Value = "RegisterClient"

$ws.Range("A242").Value = 2302
$ws.Range("B242").Value = "Client Init Date cannot be null."
$ws.Range("C242").Value = "RegisterClient"

$ws.Range("A243").Value = 2303
$ws.Range("B243").Value = "Client is already registered."
$ws.Range("C243").Value = "RegisterClient"

# --- Match the selection left behind by the original author --------------
[void]$ws.Range("C242").Select()
